$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$b1 = @{
    2 = 0.8357025263597702
    3 = 23.14408487011977
    4 = 23.42385692451
    5 = 23.16351791782924
    6 = 24.05179562076233
    7 = 26.8428576329407
    8 = 31.0010005674101
    9 = 30.30955005455831
    10 = 31.21302249326409
    11 = 32.0814784608404
    12 = 33.99042928804606
    13 = 32.69849342853592
    14 = 35.88997875391681
    15 = 36.28342292040944
    16 = 34.78004005379024
    17 = 35.06335949562434
    18 = 34.03245215152534
    19 = 34.25648546536162
    20 = 36.08369927708025
    21 = 38.1894153246502
    22 = 39.53587716389443
    23 = 40.70479306733739
    24 = 42.2264791492584
    25 = 43.09059014136039
    26 = 42.52030322815758
    27 = 41.58634909689503
    28 = 42.09063127475767
    29 = 42.56398967749087
    30 = 43.80423552521752
    31 = 43.39000644174453
    32 = 45.5711996788264
    33 = 44.56175024835782
    34 = 45.0418527143684
    35 = 46.43532306309395
    36 = 47.86006072235381
    37 = 49.23953845820868
    38 = 52.10289406769093
}
foreach ($row in $b1.Keys) {
    $ws1.Cells.Item($row, 2).Value = $b1[$row]
}

$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$b3 = @{
    2 = 0.7411965448473047
    3 = 23.04957931837449
    4 = 23.35004332539101
    5 = 23.07181194357821
    6 = 23.99071944962106
    7 = 26.81291304248416
    8 = 30.9902499966103
    9 = 30.18518440628916
    10 = 31.16386158514654
    11 = 32.00502851923844
    12 = 33.94652912936117
    13 = 32.57258373697188
    14 = 35.89177925935601
    15 = 36.17853679696357
    16 = 34.67097493986532
    17 = 34.99513835598432
    18 = 33.9208699545947
    19 = 34.18784547295445
    20 = 36.03524420792046
    21 = 38.1399798487534
    22 = 39.47005702006111
    23 = 40.64083794491979
    24 = 42.16963979578669
    25 = 43.01675190536299
    26 = 42.42085152963928
    27 = 41.48791529490769
    28 = 42.02353046785343
    29 = 42.48515933200976
    30 = 43.74643961981285
    31 = 43.28833693519783
    32 = 45.54219096167729
    33 = 44.43681444285851
    34 = 44.98356380169322
    35 = 46.37366051713227
    36 = 47.80027671472777
    37 = 49.17809475718612
    38 = 52.07473321313353
    39 = 53.20461086955618
    40 = 53.24262354462061
    41 = 54.02080619320758
    42 = 55.09758107173743
    43 = 54.97443640468325
}
foreach ($row in $b3.Keys) {
    $ws3.Cells.Item($row, 2).Value = $b3[$row]
}
